$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$prefabsView = $wb.Worksheets.Item("Prefabs View")
$prefabsView.Name = "Animals View"

$dayNightCycle = $wb.Worksheets.Item("DayNightCycle")
$dayNightCycle.Name = "Day Night Cycle"

# --- Add "Is Unique" column (H) to the "Tile Types" sheet ---
$tileTypes = $wb.Worksheets.Item("Tile Types")

$tileTypes.Range("H1").Value = "Is Unique"

$tileTypes.Range("H2").Value = $false
$tileTypes.Range("H3").Value = $false
$tileTypes.Range("H4").Value = $false
$tileTypes.Range("H5").Value = $false
$tileTypes.Range("H6").Value = $true
$tileTypes.Range("H7").Value = $true
$tileTypes.Range("H8").Value = $false

# Match the formatting of the neighbouring "Is Default" column (G) so the
# new cells carry the same cell style as the rest of the table.
$tileTypes.Range("G1:G8").Copy()
$tileTypes.Range("H1:H8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
